$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to Text format so values like "0.4660"
# or "27.716.44" are stored as literal text, not re-parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (D) updates ---
$ws.Range("D2").Value = '27.716.44'
$ws.Range("D3").Value = '1.846.89'
$ws.Range("D5").Value = '335.48'
$ws.Range("D7").Value = '0.4660'
$ws.Range("D8").Value = '0.3852'
$ws.Range("D9").Value = '46.80'
$ws.Range("D10").Value = '0.07894'
$ws.Range("D11").Value = '0.9666'
$ws.Range("D12").Value = '21.24'
$ws.Range("D13").Value = '1.848.88'
$ws.Range("D14").Value = '5.881'
$ws.Range("D15").Value = '7.133'
$ws.Range("D17").Value = '90.97'
$ws.Range("D18").Value = '0.06620'
$ws.Range("D20").Value = '17.28'
$ws.Range("D21").Value = '1.008'
$ws.Range("D22").Value = '27.727.67'
$ws.Range("D23").Value = '5.341'
$ws.Range("D25").Value = '2.296'
$ws.Range("D26").Value = '2.070.20'
$ws.Range("D27").Value = '158.79'
$ws.Range("D30").Value = '5.381'
$ws.Range("D31").Value = '118.58'
$ws.Range("D32").Value = '0.09438'
$ws.Range("D33").Value = '0.9397'
$ws.Range("D34").Value = '3.602'
$ws.Range("D35").Value = '5.247'
$ws.Range("D38").Value = '0.02210'
$ws.Range("D39").Value = '8.216'
$ws.Range("D40").Value = '1.007'
$ws.Range("D41").Value = '1.153'
$ws.Range("D42").Value = '0.5795'
$ws.Range("D43").Value = '0.1844'
$ws.Range("D45").Value = '1.284'
$ws.Range("D46").Value = '11.94'
$ws.Range("D47").Value = '0.5441'
$ws.Range("D48").Value = '1.934'
$ws.Range("D49").Value = '0.06832'
$ws.Range("D50").Value = '110.69'
$ws.Range("D51").Value = '1.008'

# --- Volume(1h) (E) updates ---
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("E7").Value = '  +0.87%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.95%  '
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("E44").Value = '  -2.75%  '
$ws.Range("E45").Value = '  +3.73%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("E51").Value = '  -32.11%  '

# --- Coin name / link swap (rows 46-47: EnergySwap <-> Decentraland) ---
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

# Restore default styling on the Price column so no stray number format
# is left applied to the cells (matches original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"

Write-Output "cryptos list updated"
